# Duplicate Sheet1 into a new "Sheet2" (placed right after Sheet1), then
# swap the contents of columns B and C on the copy, rename it, fix up its
# selection/active-cell, and finally make Sheet2 the active sheet/tab while
# Sheet1's own selection becomes the full used range.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Copy() clones formatting (row heights, styles, etc.) along with the data,
# which keeps the new sheet consistent with Sheet1 instead of picking up
# engine defaults the way Worksheets.Add() would.
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Swap columns B and C (values + number formatting) so Sheet2 shows:
#   Header 1 | Header 3 | Header 2
#   Value 11 |     1.3  |   12
#   Value 21 |     2.3  |   22
#   Value 31 |     3.3  |   32
$ws2.Range("B1").Value = "Header 3"
$ws2.Range("C1").Value = "Header 2"

$ws2.Range("B2").Value = 1.3
$ws2.Range("C2").Value = 12
$ws2.Range("B3").Value = 2.3
$ws2.Range("C3").Value = 22
$ws2.Range("B4").Value = 3.3
$ws2.Range("C4").Value = 32

$ws2.Range("B2:B4").NumberFormat = "0.00"
# B2:B4 inherited the "2 decimals" style from the copied C column; C2:C4
# inherited it the other way around, so drop C back to the default style
# (ClearFormats -> General/no explicit style, instead of materialising a
# redundant "General" numFmt).
$ws2.Range("C2:C4").ClearFormats() | Out-Null

# Sheet1's selection becomes the whole data range (no special active cell).
$ws1.Range("A1:C4").Select() | Out-Null

# Sheet2 becomes the active sheet/tab, selected at D7.
$ws2.Activate() | Out-Null
$ws2.Range("D7").Select() | Out-Null
